# "Fixed workbook and removed propane."
#
# The sheet had 6 data rows. Two of them - "Propane" (row 3) and
# "Isopentyl Acetate" (row 4) - are removed entirely (whole rows,
# shifting the rows below them up), leaving only:
#   1 Header
#   2 Carbon Dioxide
#   3 Acetaldehyde   (was row 5)
#   4 Acetic Acid    (was row 6)
#
# The header's "Paper" label in J1 is also cleared (the remaining rows'
# J-column citations stay intact), and the selected cell moves to E9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Isopentyl Acetate" row (row 4) first, then "Propane" (row 3),
# so row indices for the earlier delete aren't invalidated.
$ws.Rows("4").Delete()
$ws.Rows("3").Delete()

# The header no longer carries a "Paper" column label.
$ws.Range("J1").ClearContents()

# Restore the selection recorded in the saved workbook.
$ws.Range("E9").Select()
